$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# 1. Move the street-code lookup list from P1:P10 down to P19:P28
#    (and repoint the J3:J50 data-validation list at its new home).
# ---------------------------------------------------------------
# carry the old lookup list's look (vertical-center + wrap) down to its
# new home before wiping the old range
$ws.Range("P1:P10").Copy()
$ws.Range("P19:P28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$streetCodes = @("DNT","DTG","HL","HTLO","NT","PG","TCC","TDT","THD","XV")
for ($i = 0; $i -lt $streetCodes.Length; $i++) {
    $ws.Cells.Item(19 + $i, 16).Value = $streetCodes[$i]
}

$ws.Range("P1:P10").Clear()

$validation = $ws.Range("J3:J50").Validation
$validation.Delete()
$validation.Add(3, 1, 1, "=`$P`$19:`$P`$28")

# ---------------------------------------------------------------
# 2. Replace the sample/test data on row 3 with the real record.
# ---------------------------------------------------------------
$ws.Range("A3").Value = 123456789
$ws.Range("B3").Value = "LÊ GIA HUY"
$ws.Range("C3").Value = 25471122
$ws.Range("D3").Value = "15/12/2018"
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "bán muối"
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").Value = "HTLO"

# ---------------------------------------------------------------
# 3. Misc view state.
# ---------------------------------------------------------------
$ws.Range("J3").Select()
